# Fixed #366 User content is lost after two generation without edition.
#
# The document contains "simple field" (<w:fldSimple w:instr="...">) markers
# for the m:usercontent / m:endusercontent pseudo-fields used by M2Doc.
# Re-expand each of them into the "complex field" run sequence
# (fldChar begin / instrText / fldChar separate / fldChar end) so the field
# survives being generated a second time without the document being edited
# in between.

$d = $word.ActiveDocument

function Get-ContainingParagraph($doc, $pos) {
    for ($j = 1; $j -le $doc.Paragraphs.Count; $j++) {
        $pp = $doc.Paragraphs.Item($j)
        if ($pos -ge $pp.Range.Start -and $pos -lt $pp.Range.End) {
            return $pp
        }
    }
    return $doc.Paragraphs.Item($doc.Paragraphs.Count)
}

# Walk the fields from last to first: inserting new runs shifts character
# offsets of everything that follows, so handling later fields first keeps
# the offsets of the fields we have not processed yet valid.
for ($i = $d.Fields.Count; $i -ge 1; $i--) {
    $f = $d.Fields.Item($i)

    # Only re-expand "simple field" (wdFieldEmpty-style, i.e. fldSimple)
    # usercontent markers; leave any other kind of field alone.
    $instr = $f.Code.Text
    if ($instr -ne "m:usercontent zone1" -and $instr -ne "m:endusercontent") {
        continue
    }

    $fieldPara = Get-ContainingParagraph $d $f.Code.Start
    $insertPoint = $d.Range($fieldPara.Range.Start, $fieldPara.Range.Start)

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
           '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' +
           '<w:r><w:instrText>' + $instr + '</w:instrText></w:r>' +
           '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' +
           '<w:r><w:fldChar w:fldCharType="end"/></w:r>' +
           '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $insertPoint.InsertXML($xml)

    # The original fldSimple is still present right after the freshly
    # inserted runs (it is now the next field in the collection) -- remove
    # it, leaving only the new begin/instrText/separate/end runs behind.
    $stale = $d.Fields.Item($i + 1)
    $stale.Delete()
}
